$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8

# Row 3
$ws.Range("G3").Value = 2.05
$ws.Range("I3").Value = 3.3
$ws.Range("J3").Value = 2.75
$ws.Range("L3").Value = 4
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 3.4
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.85
$ws.Range("X3").Value = 10
$ws.Range("Z3").Value = 19
$ws.Range("AC3").Value = 10
$ws.Range("AE3").Value = 15
$ws.Range("AF3").Value = 51
$ws.Range("AI3").Value = 12
$ws.Range("AK3").Value = 26
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 11
$ws.Range("AU3").Value = 8
$ws.Range("AW3").Value = 5.5
$ws.Range("AX3").Value = 19
$ws.Range("BA3").Value = 81
$ws.Range("BB3").Value = 201

# Row 4
$ws.Range("G4").Value = 2.62
$ws.Range("H4").Value = 2.8
$ws.Range("I4").Value = 2.75
$ws.Range("J4").Value = 3.25
$ws.Range("K4").Value = 1.98
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.72
$ws.Range("Q4").Value = 2.22
$ws.Range("R4").Value = 1.62
$ws.Range("T4").Value = 2.55
$ws.Range("U4").Value = 1.82
$ws.Range("W4").Value = 7.5
$ws.Range("X4").Value = 13
$ws.Range("Y4").Value = 9.75
$ws.Range("Z4").Value = 32
$ws.Range("AA4").Value = 24
$ws.Range("AC4").Value = 6
$ws.Range("AD4").Value = 5.5
$ws.Range("AF4").Value = 70
$ws.Range("AG4").Value = 7.6
$ws.Range("AH4").Value = 13.5
$ws.Range("AL4").Value = 37
$ws.Range("AN4").Value = 4.55
$ws.Range("AO4").Value = 14.5
$ws.Range("AQ4").Value = 65
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.55
$ws.Range("AU4").Value = 6.7
$ws.Range("AV4").Value = 60
$ws.Range("AW4").Value = 4.7

# Row 7
$ws.Range("G7").Value = 5.6
$ws.Range("J7").Value = 5.7
$ws.Range("K7").Value = 2.22
$ws.Range("L7").Value = 2.07
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 7.1
$ws.Range("O7").Value = 1.32
$ws.Range("P7").Value = 3.1
$ws.Range("Q7").Value = 1.98
$ws.Range("R7").Value = 1.78
$ws.Range("S7").Value = 1.38
$ws.Range("T7").Value = 2.8
$ws.Range("U7").Value = 2.02
$ws.Range("W7").Value = 14
$ws.Range("AC7").Value = 7.1
$ws.Range("AE7").Value = 19
$ws.Range("AG7").Value = 5.9
$ws.Range("AH7").Value = 6.6
$ws.Range("AJ7").Value = 10.75
$ws.Range("AK7").Value = 13.5
$ws.Range("AP7").Value = 37
$ws.Range("AQ7").Value = 200
$ws.Range("AT7").Value = 2.8
$ws.Range("AU7").Value = 8
$ws.Range("AV7").Value = 75
$ws.Range("AX7").Value = 7.3
$ws.Range("AY7").Value = 17.5
$ws.Range("AZ7").Value = 23
$ws.Range("BA7").Value = 55
$ws.Range("BB7").Value = 250
